# Generate Report for Handoff
#
# The "0686809e-566d-475d-8c26-c940e58fc9a0" localization file has finished
# its en-US sync and is now ready to be handed off again, while the
# "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c" file remains handed-back / in sync.
# The status report re-sorts these two file rows (f01f6e9d now listed
# before 0686809e on every sheet) and updates 0686809e's status + handoff
# timestamp accordingly.

$wb = $excel.ActiveWorkbook

$file_0686 = "0686809e-566d-475d-8c26-c940e58fc9a0.md"
$file_f01f = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.md"

function Set-HyperlinkDisplay($ws, $addr, $text) {
    $map = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $map[$hl.Range.Address(0, 0)] = $hl
    }
    $map[$addr].TextToDisplay = $text
}

# ---------------------------------------------------------------------
# Overview sheet: File Name / zh-cn / de-de summary
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = $file_f01f
$ov.Range("B2").Value = "Handed back: in sync with en-US"
$ov.Range("C2").Value = "Handed back: in sync with en-US"

$ov.Range("A3").Value = $file_0686
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

Set-HyperlinkDisplay $ov "A2" $file_f01f
Set-HyperlinkDisplay $ov "A3" $file_0686

# ---------------------------------------------------------------------
# zh-cn sheet: per-language handoff/handback detail
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = $file_f01f
$zh.Range("B2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf"
$zh.Range("D2").Value = "2016-03-10 14:37:32"
$zh.Range("E2").Value = $file_f01f
$zh.Range("F2").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf"
$zh.Range("G2").Value = "2016-03-10 14:39:16"
$zh.Range("H2").Value = "Include"

$zh.Range("A3").Value = $file_0686
$zh.Range("B3").Value = "Ready for handoff"
$zh.Range("C3").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf"
$zh.Range("D3").Value = "2016-03-10 14:39:45"
$zh.Range("E3").Value = $file_0686
$zh.Range("F3").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf"
$zh.Range("G3").Value = "2016-03-10 14:39:16"
$zh.Range("H3").Value = "Include"

Set-HyperlinkDisplay $zh "A2" $file_f01f
Set-HyperlinkDisplay $zh "C2" "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf"
Set-HyperlinkDisplay $zh "E2" $file_f01f
Set-HyperlinkDisplay $zh "F2" "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.zh-cn.xlf"

Set-HyperlinkDisplay $zh "A3" $file_0686
Set-HyperlinkDisplay $zh "C3" "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf"
Set-HyperlinkDisplay $zh "E3" $file_0686
Set-HyperlinkDisplay $zh "F3" "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.zh-cn.xlf"

# ---------------------------------------------------------------------
# de-de sheet: per-language handoff/handback detail
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = $file_f01f
$de.Range("B2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf"
$de.Range("D2").Value = "2016-03-10 14:38:16"
$de.Range("E2").Value = $file_f01f
$de.Range("F2").Value = "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf"
$de.Range("G2").Value = "2016-03-10 14:39:23"
$de.Range("H2").Value = "Include"

$de.Range("A3").Value = $file_0686
$de.Range("B3").Value = "Ready for handoff"
$de.Range("C3").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf"
$de.Range("D3").Value = "2016-03-10 14:39:49"
$de.Range("E3").Value = $file_0686
$de.Range("F3").Value = "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf"
$de.Range("G3").Value = "2016-03-10 14:39:23"
$de.Range("H3").Value = "Include"

Set-HyperlinkDisplay $de "A2" $file_f01f
Set-HyperlinkDisplay $de "C2" "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf"
Set-HyperlinkDisplay $de "E2" $file_f01f
Set-HyperlinkDisplay $de "F2" "f01f6e9d-34cd-4a5f-9bfe-4aba36a9b13c.e6b5ef5ec3c4c2ff412ccabdad540ba8efea9d84.de-de.xlf"

Set-HyperlinkDisplay $de "A3" $file_0686
Set-HyperlinkDisplay $de "C3" "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf"
Set-HyperlinkDisplay $de "E3" $file_0686
Set-HyperlinkDisplay $de "F3" "0686809e-566d-475d-8c26-c940e58fc9a0.c33c22caa5a0c1e9a12e1d808322b661a4f4e7f0.de-de.xlf"
